$ns = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'
$d = $word.ActiveDocument

function Set-ParagraphXml {
    param(
        [int]$Index,
        [string]$ExpectedSubstring,
        [string]$InnerXml
    )
    $para = $d.Paragraphs($Index)
    $rng = $para.Range
    if ($ExpectedSubstring -ne "" -and $rng.Text.IndexOf($ExpectedSubstring) -lt 0) {
        throw "Paragraph $Index does not contain expected text '$ExpectedSubstring' (got: $($rng.Text))"
    }
    $xml = '<w:p xmlns:w="' + $ns + '">' + $InnerXml + '</w:p>'
    $rng.InsertXML($xml) | Out-Null
}

# Paragraph 6: "Create the correct xpath..." - re-split runs, switch most of the
# line to en-US language, and add spellStart/spellEnd + gramStart/gramEnd proofing marks.
$p6Inner = '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>' + [char]0x421 + [char]0x43E + [char]0x437 + [char]0x434 + [char]0x430 + [char]0x43B + [char]0x438 + '</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>' + [char]0x441 + [char]0x43B + [char]0x435 + [char]0x434 + [char]0x443 + [char]0x44E + [char]0x449 + [char]0x438 + [char]0x439 + '</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>' + [char]0x43F + [char]0x440 + [char]0x430 + [char]0x432 + [char]0x438 + [char]0x43B + [char]0x44C + [char]0x43D + [char]0x44B + [char]0x439 + '</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>xpath</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>: .</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>//</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>div</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>[@</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>class</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>=' + [char]0x27 + '</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>page</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>-</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>artist</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>' + [char]0x27 + ']//</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>div</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>[@</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>class</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>=' + [char]0x27 + '</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>album</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>album</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>_</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>selectable</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>' + [char]0x27 + ']//</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>img</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>[@</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>class</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>=' + [char]0x27 + '</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>album</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>-</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>cover</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>album</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>-</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>cover</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>_</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>size</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>_</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>L</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>' + [char]0x27 + ']</w:t></w:r>'

Set-ParagraphXml 6 "xpath" $p6Inner

# Paragraph 7: "Div - type of element in the tree." - wrap "Div" with spellStart/spellEnd.
$p7Inner = '<w:pPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Div</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> ' + [char]0x2013 + ' </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>' + [char]0x442 + [char]0x438 + [char]0x43F + ' ' + [char]0x44D + [char]0x43B + [char]0x435 + [char]0x43C + [char]0x435 + [char]0x43D + [char]0x442 + [char]0x430 + ' ' + [char]0x432 + ' ' + [char]0x434 + [char]0x435 + [char]0x440 + [char]0x435 + [char]0x432 + [char]0x435 + '.</w:t></w:r>'

Set-ParagraphXml 7 "Div" $p7Inner

# Paragraph 10: "Shortened xpath: (.//img)[4]" - re-split runs, switch paragraph
# mark language back to ru-RU, and add proofing marks.
$p10Inner = '<w:pPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">' + [char]0x421 + [char]0x43E + [char]0x43A + [char]0x440 + [char]0x430 + [char]0x449 + [char]0x435 + [char]0x43D + [char]0x43D + [char]0x44B + [char]0x439 + ' </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>xpath</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>(./</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>/</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>img</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>)[4]</w:t></w:r>'

Set-ParagraphXml 10 "xpath" $p10Inner

# Drop the trailing empty paragraph (the one right after the _GoBack bookmark
# paragraph, just before the section break). Its own mark cannot be deleted in
# isolation, so extend the deleted range back into the previous paragraph mark.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
if ($lastPara.Range.Text.Trim() -eq "") {
    $prevPara = $d.Paragraphs($count - 1)
    $delRange = $d.Range($prevPara.Range.End - 1, $lastPara.Range.End)
    $delRange.Delete() | Out-Null
}
